$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.49 = 25554.81 pesos`n✅ 25554.81 pesos = 6.44 = 940.32 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 154
$wsTasas.Range("O10").Value = 3935.44
$wsTasas.Range("N12").Value = 3967.8
$wsTasas.Range("O12").Value = 146
